$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results (case with 380 kV) for rows 2-25, columns B-O (excluding E,L,N which stay 0, and A which is the index)
$data = @{
    2 = @(0.273140775175392,0.04209342303292374,0.3221304964359035,1.195002854771417,0.5616081605313994,0.7071845117177347,0.5692729848684088,0.3157874387711033,0.2824326445084466,0.2547900707812261,2.51270601604746)
    3 = @(0.2394922067879293,0.03749247042848936,0.3161389828031247,1.200611460957106,0.5663024528189169,0.712884750484001,0.5758924875430331,0.314596671398462,0.2465844776632196,0.2401114562535795,2.534359159853665)
    4 = @(0.2187791219357393,0.03465198786713586,0.3125925338377442,1.204679329062934,0.5695445945006554,0.7166684454153085,0.5802525382093613,0.3140312065565283,0.2244974574680185,0.2311744579039967,2.54900304331909)
    5 = @(0.2103257061969828,0.03349065121301464,0.3111807643275029,1.206494078425024,0.5709562451729226,0.7182817370842045,0.5821036359288438,0.3138424719159119,0.2154782962732469,0.2275518532978182,2.555309760171014)
    6 = @(0.2089212754808614,0.03329758383468118,0.3109483651041387,1.206804905436805,0.5711961108665449,0.7185539374111372,0.5824154992984898,0.3138136521646189,0.21397957093005,0.2269514953324432,2.556377475163941)
    7 = @(0.2186651667089166,0.03463634102976698,0.3125733586310702,1.204703167331786,0.5695632663201806,0.7166899136335516,0.5802772017561644,0.3140284923363339,0.2243758960812698,0.2311255237334464,2.549086724362283)
    8 = @(0.2615500378285844,0.04051027106308425,0.3200372396219251,1.196807221917133,0.5631520521661102,0.7090911031577818,0.5714940268073914,0.3153424985053661,0.2700883490935837,0.2497133128275877,2.519892117856983)
    9 = @(0.3452066854587201,0.05190342717902752,0.3357184321364599,1.186271691667102,0.5534363434821614,0.6964384266894186,0.5566160707170802,0.3192326689389233,0.3591036876878775,0.2867556773150071,2.473341482619773)
    10 = @(0.4063768984965179,0.06019444524160633,0.347870010266405,1.181543620308112,0.5480422767977942,0.6885096827205501,0.5471151432949632,0.3228907486848271,0.4240967602832768,0.314322064675828,2.445661905351741)
    11 = @(0.4341366553355215,0.0639483953354727,0.3535337855445562,1.180045975802706,0.5459676569511558,0.6851988060662748,0.5431033812695851,0.3247285054748659,0.4535707235451696,0.3269372001635134,2.434485468449509)
    12 = @(0.4446384115598789,0.06536731083677694,0.3556979252563792,1.179572695522047,0.545236618547861,0.6839875637575901,0.5416288515752861,0.3254493654529682,0.4647180168144303,0.331724797074024,2.430456714959618)
    13 = @(0.4423771370619534,0.06506183999974269,0.3552309791843982,1.17967045241074,0.5453916326994559,0.6842465362730579,0.5419444322996263,0.3252930066201998,0.4623178746813608,0.3306932389344155,2.431315327091966)
    14 = @(0.4350008504736991,0.06406518342424761,0.3537114429163069,1.180005158419874,0.5459064199248829,0.685098304613156,0.5429811760913701,0.3247873113893576,0.4544880996103302,0.3273308693920995,2.434149941481067)
    15 = @(0.430481306554185,0.06345435834725777,0.3527832035900929,1.180222394504938,0.5462288507599951,0.6856255727506806,0.5436220253979123,0.3244808055650168,0.4496903118928799,0.3252726827294126,2.435912729749418)
    16 = @(0.4045613322189183,0.05994875306224401,0.3475025914404171,1.181654631613,0.5481854923110063,0.6887320079259283,0.5473835665102555,0.3227741389631547,0.4221686647715899,0.3134991227309953,2.44642079022077)
    17 = @(0.3886426792169289,0.05779359204592538,0.3442978160539099,1.182700500030492,0.5494829792718932,0.6907134747403632,0.5497706308506629,0.3217716138060212,0.4052610596166915,0.3062954524036599,2.453229614883398)
    18 = @(0.3794804343742442,0.05655234161389444,0.3424673209429585,1.18336354107737,0.5502649455461679,0.6918810206529002,0.5511728086705787,0.3212113372625538,0.3955276598718171,0.3021591719187526,2.457279069162553)
    19 = @(0.3763772007617945,0.05613179315139405,0.3418497512280254,1.183598597982808,0.5505358334388077,0.6922811176714276,0.5516525766830753,0.321024446287467,0.3922306426139528,0.3007599227534641,2.458673020979049)
    20 = @(0.3903378992160356,0.0580231848928463,0.3446376455550961,1.182582802649819,0.5493411658197473,0.6904996612608159,0.5495135013599501,0.3218766424579087,0.4070617971865715,0.3070615646144432,2.452491018151775)
    21 = @(0.4371677282590269,0.0643579974874342,0.3541572425666288,1.179904300871442,0.5457537328817565,0.6848469660996983,0.5426754476777731,0.3249351697028828,0.4567882757406778,0.3283181955275722,2.433311822840068)
    22 = @(0.4677136207582748,0.06848283151263956,0.360491800180796,1.178700708393251,0.5437272667282471,0.6814003925160321,0.5384665764288954,0.3270794556745926,0.4892063260320469,0.3422718180112909,2.421963394419024)
    23 = @(0.4514164152726607,0.06628275994772537,0.3571006455687211,1.179293068224197,0.5447797043624831,0.6832172348548013,0.5406891151098066,0.3259217204156641,0.4719118416794288,0.3348189988156136,2.427911711256343)
    24 = @(0.3895715225977767,0.05791939287330194,0.3444839713321954,1.182635821260668,0.5494051674741769,0.6905962379479718,0.5496296567088628,0.3218291088881813,0.4062477241504325,0.3067151890493065,2.4528245173139)
    25 = @(0.3226249681899276,0.04883503064992567,0.3313650659642491,1.188592513024595,0.555758606921998,0.6996209604297405,0.5603898329513086,0.318039738207176,0.3350923046241689,0.2766723545640133,2.484789344755953)
}

$cols = @("B","C","D","F","G","H","I","J","K","M","O")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

Write-Host "Updated $($data.Keys.Count) rows"